$wb = $excel.ActiveWorkbook

$nodes = $wb.Worksheets.Item("nodes")
$edges = $wb.Worksheets.Item("edges")

# "added bifurcation option to main script" -> new "fraction" column (D)
# on the edges sheet, populated with 1 for every existing edge row (2-10)
for ($r = 2; $r -le 10; $r++) {
    $edges.Cells.Item($r, 4).Value = 1
}

# Keep the nodes-sheet selection as-is (E12)
$nodes.Range("E12").Select() | Out-Null

# Switch the active/selected sheet to "edges" (tabSelected + activeTab)
# and move its selection to D11
$edges.Activate()
$edges.Range("D11").Select() | Out-Null
